# Update "想去人数" (want-to-go count) values in the "F" column
# for both the "展览" and "全部类型" worksheets, which hold the same data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of cell address -> new value
$updates = @{
    "F3"  = 2634
    "F4"  = 533
    "F6"  = 6598
    "F7"  = 446
    "F11" = 4
    "F12" = 133
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
